$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 20; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value2 = 45183
}
